# Add 6 new location rows (236-241) to the "location-1" sheet, matching the
# "Add files via upload" commit that appended more LIVE-cam entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Category, B=latitude,longitude, C=Location, D=CITY, E=COUNTRY, F=YouTube Link

$newRows = @(
    @{ A = "LIVE, CITY, BUILDING";
       B = "37.55131141702741, 126.98834538813063";
       C = "남산서울타워 파노라마 LIVE / Namsan Seoul Tower Panorama";
       D = "Seoul";
       E = "South Korea";
       F = "CegCJLdx3aw" },
    @{ A = "LIVE, RELIGION";
       B = "29.99230752240382, 78.19167894946095";
       C = "लाइव दर्शन || गायत्री तीर्थ शांतिकुंज हरिद्वार / Gayatri Teerth Shantikunj Haridwar";
       D = "Uttarakhand";
       E = "India";
       F = "6pKI_gV2ibQ" },
    @{ A = "LIVE, RELIGION";
       B = "25.310865117268058, 83.01067751835699";
       C = "🔴Live Darshan From Shree Kashi Vishwanath Temple Varanasi ( श्री काशी विश्वनाथ मंदिर से लाइव दर्शन )";
       D = "Uttar Pradesh";
       E = "India";
       F = "jI_kv9h2sZM" },
    @{ A = "LIVE, RELIGION";
       B = "25.596062386202732, 85.22995712137404";
       C = "🔴 Takht Sri Patna Sahib ( Live Darshan from Harmandir ) (गुरुद्वारा श्री हरिमंदर जी पटना साहिब) !!!";
       D = "Bihar";
       E = "India";
       F = "6ON11RBhWME" },
    @{ A = "LIVE, RELIGION";
       B = "26.92885651857421, 75.82407006189108";
       C = "🔴 Thikana Mandir Shri GOVIND DEVJI, Jaipur LIVE DARSHAN";
       D = "Rajasthan";
       E = "India";
       F = "_9jLnXrRDUw" },
    @{ A = "LIVE, RELIGION";
       B = "23.182902707115378, 75.76824133445302";
       C = "🔴Live Darshan - Shree Mahakaleshwar Temple Ujjain (महाकालेश्वर मंदिर के लाइव दर्शन) !";
       D = "Madhya Pradesh";
       E = "India";
       F = "eBikK8yBOtI" }
)

$startRow = 236
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Carry the existing left/right-border cell style down from the row above
    # (used throughout the table) before filling in the new values.
    $ws.Range("A" + ($r - 1) + ":F" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = $row.A
    $ws.Range("B" + $r).Value = $row.B
    $ws.Range("C" + $r).Value = $row.C
    $ws.Range("D" + $r).Value = $row.D
    $ws.Range("E" + $r).Value = $row.E
    $ws.Range("F" + $r).Value = $row.F
}

$excel.CutCopyMode = 0

# Match the final saved selection from the source workbook.
$ws.Range("F243").Select()
